$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 15
$ws.Range("C3").Value = 9
$ws.Range("C5").Value = 17
$ws.Range("C6").Value = 15
$ws.Range("C7").Value = 19
$ws.Range("C8").Value = 15
$ws.Range("C10").Value = 15
$ws.Range("C11").Value = 9
$ws.Range("C12").Value = 6
$ws.Range("C13").Value = 13
$ws.Range("C14").Value = 12
$ws.Range("C15").Value = 16
$ws.Range("C16").Value = 12
$ws.Range("C17").Value = 19
$ws.Range("B18").Value = "<in>"
